$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two shared-string values used by B1 and C1
$ws.Range("B1").Value = "-МЕНЮ"
$ws.Range("C1").Value = "Основное меню=yes"

# Move the active selection from B2 to C2
$ws.Range("C2").Select()
